# Update "想去人数" (wanted-to-go count) figures in the three data sheets
# (展览, 演出, 全部类型) to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibition) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 33
$ws1.Range("F4").Value  = 1429
$ws1.Range("F5").Value  = 175
$ws1.Range("F6").Value  = 33
$ws1.Range("F7").Value  = 33
$ws1.Range("F8").Value  = 9616
$ws1.Range("F9").Value  = 161
$ws1.Range("F11").Value = 238
$ws1.Range("F13").Value = 364
$ws1.Range("F14").Value = 6614
$ws1.Range("F15").Value = 1078
$ws1.Range("F16").Value = 119
$ws1.Range("F18").Value = 173

# ---- Sheet: 演出 (Performance) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 42

# ---- Sheet: 全部类型 (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 33
$ws4.Range("F4").Value  = 1429
$ws4.Range("F5").Value  = 175
$ws4.Range("F6").Value  = 33
$ws4.Range("F7").Value  = 33
$ws4.Range("F8").Value  = 42
$ws4.Range("F10").Value = 9616
$ws4.Range("F11").Value = 161
$ws4.Range("F13").Value = 238
$ws4.Range("F15").Value = 364
$ws4.Range("F16").Value = 6614
$ws4.Range("F17").Value = 1078
$ws4.Range("F18").Value = 119
$ws4.Range("F20").Value = 173
